# Add 8 new "grocery" XPath rows (rows 72-79) to the "XPath" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("XPath")

$rows = @(
  @("grocery active count", '//*[@id="root"]/div/div[4]/div[1]/div/div/div/div[2]/div[1]/div[3]/div[2]/div/div[2]/div/div/div[2]/div[1]/span[1]', 12.8),
  @("grocery delivered count", '//*[@id="root"]/div/div[4]/div[1]/div/div/div/div[2]/div[1]/div[3]/div[2]/div/div[2]/div/div/div[2]/div[1]/span[2]', 23.05),
  @("grocery cancelled count", '//*[@id="root"]/div/div[4]/div[1]/div/div/div/div[2]/div[1]/div[3]/div[2]/div/div[2]/div/div/div[2]/div[1]/span[3]', 23.05),
  @("grocery returned count", '//*[@id="root"]/div/div[4]/div[1]/div/div/div/div[2]/div[1]/div[3]/div[2]/div/div[2]/div/div/div[2]/div[1]/span[4]', 23.05),
  @("grocery undelivered count", '//*[@id="root"]/div/div[4]/div[1]/div/div/div/div[2]/div[1]/div[3]/div[2]/div/div[2]/div/div/div[2]/div[1]/span[5]', 23.05),
  @("grocery basket price", '//*[@id="root"]/div/div[4]/div[1]/div/div/div/div[2]/div[1]/div[3]/div[2]/div/div[2]/div/div/div[2]/div[2]/div/div[1]', 12.8),
  @("grocery delivery charge", '//*[@id="root"]/div/div[4]/div[1]/div/div/div/div[2]/div[1]/div[3]/div[2]/div/div[2]/div/div/div[2]/div[2]/div/div[2]', 23.05),
  @("grocery total saving", '//*[@id="root"]/div/div[4]/div[1]/div/div/div/div[2]/div[1]/div[3]/div[2]/div/div[2]/div/div/div[2]/div[2]/div/div[3]', 12.8)
)

$r = 72
foreach ($item in $rows) {
  $ws.Cells.Item($r, 1).Value = $item[0]
  $ws.Cells.Item($r, 2).Value = $item[1]
  $ws.Rows.Item($r).RowHeight = $item[2]
  $r = $r + 1
}

$ws.Range("B80").Select()
